# Remove the "divisi" and "posisi" columns (O:P) from Sheet1 and shift the
# remaining columns (sisa_cuti, status, tanggal_gabung, ...) left by two.
# This matches the diff: the two <si> entries "divisi"/"posisi" (and their
# sample data "Board Of Director"/"Software Engineer") disappear from
# sharedStrings.xml, every column from Q onward shifts to O onward, the
# sheet dimension shrinks from A1:W2 to A1:U2, and the two data validations
# that used to target Q2/R2 now target O2/P2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Columns("O:P").Delete()

# Restore the active selection recorded in the saved file.
[void]$ws.Range("G10").Select()
